$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(2)

# Force a genuine content change so the save/diff logic actually rewrites
# the run structure (merging the 3 existing runs into 1), instead of
# treating an identical concatenated-text assignment as a no-op.
$para.Text = "__tmp__"
$para.Text = "Actually, this project will be get popularity who live in urban area."
